# 001 H - O sebi.pptx small edits
#  1. Refresh the cached "Update automatically" date/time field text
#     (datetimeFigureOut) from 03/05/2014 to 08/05/2014 on the slide
#     master and every slide layout's Date placeholder.
#  2. Type the subtitle "tko je taj tip?" on slide 1.
#  3. Split "Who's that guy" into two runs on slide 2's title.

$p = $ppt.ActivePresentation

function Update-DateTimePlaceholder {
    param($shapes)

    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        $isDate = $false
        try {
            if ($shp.PlaceholderFormat.Type -eq 16) {
                $isDate = $true
            }
        } catch {
        }
        if ($isDate) {
            $shp.TextFrame.TextRange.Text = "08/05/2014"
        }
    }
}

# --- 1. Date placeholder on the master and on every layout ---
Update-DateTimePlaceholder $p.SlideMaster.Shapes

for ($li = 1; $li -le $p.SlideMaster.CustomLayouts.Count; $li++) {
    $layout = $p.SlideMaster.CustomLayouts.Item($li)
    Update-DateTimePlaceholder $layout.Shapes
}

# --- 2. Slide 1 (title "Hrvoje") - type the subtitle text ---
$slide1 = $p.Slides.Item(1)
$subtitle = $slide1.Shapes.Item(2)
$subtitleRange = $subtitle.TextFrame.TextRange
$subtitleRange.InsertAfter("tko je taj tip?") | Out-Null

$w1 = $subtitleRange.Characters(1, 3)
$w1.LanguageID = "en-GB"
$w1.Text = "tko"

$w2 = $subtitleRange.Characters(4, 4)
$w2.LanguageID = "en-GB"
$w2.Text = " je "

$w3 = $subtitleRange.Characters(8, 3)
$w3.LanguageID = "en-GB"
$w3.Text = "taj"

$w4 = $subtitleRange.Characters(11, 5)
$w4.LanguageID = "en-GB"
$w4.Text = " tip?"

# --- 3. Slide 2 (title "Who's that guy") - split into two runs ---
$slide2 = $p.Slides.Item(2)
$title2 = $slide2.Shapes.Item(1)
$titleRange = $title2.TextFrame.TextRange
$len = $titleRange.Length
$lastWord = $titleRange.Characters($len - 2, 3)
$lastWord.Text = "guy"
